{"js": "// 1) Fix the run text \"0-\" -> \"-\" (so \"70-minute\" becomes \"7-minute\").\nconst body = context.document.body;\nconst zeroDashResults = body.search(\"0-\", { matchCase: true });\nzeroDashResults.load(\"items\");\nawait context.sync();\nif (zeroDashResults.items.length > 0) {\n  zeroDashResults.items[0].insertText(\"-\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Move the \"_GoBack\" bookmark from that spot to span from the start of\n//    the \"Rubric:\" paragraph through the end of the \"...very strict on this\n//    point.\" paragraph.\ntry {\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n} catch (e) {\n  // no-op if it doesn't currently exist\n}\n\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items.forEach((p) => p.load(\"text\"));\nawait context.sync();\n\nlet rubricParagraph = null;\nlet strictParagraph = null;\nfor (const p of paragraphs.items) {\n  if (rubricParagraph === null && p.text.indexOf(\"Rubric:\") >= 0) {\n    rubricParagraph = p;\n  }\n  if (p.text.indexOf(\"very strict on this point\") >= 0) {\n    strictParagraph = p;\n  }\n}\n\nif (rubricParagraph && strictParagraph) {\n  const startRange = rubricParagraph.getRange(\"Start\");\n  const endRange = strictParagraph.getRange(\"End\");\n  const bookmarkRange = startRange.expandTo(endRange);\n  bookmarkRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Fix the run text \"0-\" -> \"-\" (so \"70-minute\" becomes \"7-minute\").\n$fixRange = $d.Content.Duplicate\n$fixRange.Find.Execute(\"0-minute\") | Out-Null\nif ($fixRange.Find.Found) {\n    # Only replace the \"0-\" part, leave \"minute\" untouched.\n    $zeroDashRange = $d.Range($fixRange.Start, $fixRange.Start + 2)\n    $zeroDashRange.Text = \"-\"\n}\n\n# 2) Move the \"_GoBack\" bookmark from that spot to span from the start of\n#    the \"Rubric:\" paragraph through the end of the \"...very strict on this\n#    point.\" paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$rubricRange = $d.Content.Duplicate\n$rubricRange.Find.Execute(\"Rubric:\") | Out-Null\n$startPos = $rubricRange.Start\n\n$strictRange = $d.Content.Duplicate\n$strictRange.Find.Execute(\"very strict on this point\") | Out-Null\n$endPos = $strictRange.End\n\n$bookmarkRange = $d.Range($startPos, $endPos)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
